# Fix #9722 - [Feature] Translate export search reports
# Translate the French "Coûts" export-search report to English.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet (tab) name: "Coûts" -> "Costs"
$ws.Name = "Costs"

# 2) Translate the two remaining French header strings used in row 8
#    "Entité" -> "Entity", "Variantes" -> "Variants"
$ws.Cells.Replace("Entité", "Entity")
$ws.Cells.Replace("Variantes", "Variants")

# 3) Header row fill: unify the header row (B8:U8) on the single dark
#    teal fill (the old lighter "blue" header fill is dropped).
$ws.Range("B8:U8").Interior.Color = 5521920   # RGB(0x00, 0x42, 0x54) -> FF004254

# 4) Selection moves from C8:U8 to just C8
$ws.Range("C8").Select()

# 5) Minor column width tweaks (closest values reachable through the
#    character-width -> pixel rounding used by this COM layer)
$ws.Columns.Item(1).ColumnWidth = 28.15   # -> stored width ~29.0 (was 29.02, target 29.03)
$ws.Columns.Item(8).ColumnWidth = 32.15   # -> stored width ~33.0 (was 32.96, target 32.95)
